# Update "想去人数" (interest count, column F) values on the "展览" and
# "全部类型" sheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 108
$ws1.Range("F3").Value = 12099
$ws1.Range("F4").Value = 47
$ws1.Range("F7").Value = 229
$ws1.Range("F8").Value = 11997
$ws1.Range("F9").Value = 504
$ws1.Range("F12").Value = 590
$ws1.Range("F13").Value = 1801
$ws1.Range("F14").Value = 5931

# --- Sheet "全部类型" (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 108
$ws4.Range("F5").Value = 12099
$ws4.Range("F6").Value = 47
$ws4.Range("F10").Value = 229
$ws4.Range("F11").Value = 11997
$ws4.Range("F12").Value = 504
$ws4.Range("F15").Value = 590
$ws4.Range("F16").Value = 1801
$ws4.Range("F18").Value = 5931
